$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 829 (2026/12/29 ...), shifting
# all subsequent rows down by one.
$ws.Rows.Item(829).Insert()

# The newly inserted row needs the date/day-of-week columns stored as
# literal text (matching the rest of the sheet) rather than being
# auto-converted to a date serial number, so force a text number format
# before writing the values.
$ws.Cells.Item(829, 1).NumberFormat = "@"
$ws.Cells.Item(829, 2).NumberFormat = "@"

$ws.Cells.Item(829, 1).Value = "2026/02/16"
$ws.Cells.Item(829, 2).Value = "月"
$ws.Cells.Item(829, 3).Value = 7
$ws.Cells.Item(829, 4).Value = 201
